# Apply the cryptos-list price/volume refresh described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.500.12"
$ws.Range("E2").Value = "  +3.16%  "
$ws.Range("D3").Value = "2.301.96"
$ws.Range("E3").Value = "  +2.27%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'510.33"
$ws.Range("E5").Value = "  +3.63%  "
$ws.Range("D6").Value = "'130.39"
$ws.Range("E6").Value = "  +2.33%  "
$ws.Range("D7").Value = "'0.993"
$ws.Range("E7").Value = "  -0.66%  "
$ws.Range("D8").Value = "'0.531"
$ws.Range("E8").Value = "  +1.72%  "
$ws.Range("D9").Value = "2.331.50"
$ws.Range("E9").Value = "  +3.69%  "
$ws.Range("D10").Value = "'0.101"
$ws.Range("E10").Value = "  +7.65%  "
$ws.Range("E11").Value = "  +1.27%  "
$ws.Range("D12").Value = "'5.12"
$ws.Range("E12").Value = "  +8.59%  "
$ws.Range("D13").Value = "'0.343"
$ws.Range("E13").Value = "  +3.17%  "
$ws.Range("D14").Value = "'23.81"
$ws.Range("E14").Value = "  +6.10%  "
$ws.Range("D15").Value = "2.706.14"
$ws.Range("E15").Value = "  +2.31%  "
$ws.Range("D16").Value = "55.405.25"
$ws.Range("E16").Value = "  +3.04%  "
$ws.Range("D17").Value = "'0.0000133"
$ws.Range("E17").Value = "  +3.13%  "
$ws.Range("D18").Value = "2.338.69"
$ws.Range("E18").Value = "  +3.60%  "
$ws.Range("D19").Value = "'10.47"
$ws.Range("E19").Value = "  +3.17%  "
$ws.Range("D20").Value = "'4.22"
$ws.Range("E20").Value = "  +2.72%  "
$ws.Range("D21").Value = "'317.59"
$ws.Range("E21").Value = "  +5.48%  "
$ws.Range("D22").Value = "'6.59"
$ws.Range("E22").Value = "  +5.36%  "
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("D24").Value = "'60.30"
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("D25").Value = "'0.990"
$ws.Range("E25").Value = "  -1.63%  "
$ws.Range("D26").Value = "'0.157"
$ws.Range("E26").Value = "  +6.80%  "
$ws.Range("D27").Value = "'7.61"
$ws.Range("E27").Value = "  +5.23%  "
$ws.Range("D28").Value = "'171.52"
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("D29").Value = "'6.20"
$ws.Range("E29").Value = "  +5.71%  "
$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").Value = "'1.17"
$ws.Range("E30").Value = "  +9.85%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.66"
$ws.Range("E31").Value = "  +4.56%  "
$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").Value = "0.0₃0717"
$ws.Range("E32").Value = "  +5.50%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'18.12"
$ws.Range("E33").Value = "  +2.60%  "
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "'0.992"
$ws.Range("E35").Value = "  -0.60%  "
$ws.Range("D36").Value = "'1.25"
$ws.Range("E36").Value = "  +6.18%  "
$ws.Range("D37").Value = "'0.915"
$ws.Range("E37").Value = "  -2.03%  "
$ws.Range("D38").Value = "'3.94"
$ws.Range("E38").Value = "  +7.43%  "
$ws.Range("D39").Value = "'37.08"
$ws.Range("E39").Value = "  +3.47%  "
$ws.Range("D40").Value = "'1.48"
$ws.Range("E40").Value = "  +7.48%  "
$ws.Range("D41").Value = "'0.378"
$ws.Range("E41").Value = "  +2.51%  "
$ws.Range("D42").Value = "'137.27"
$ws.Range("E42").Value = "  +10.78%  "
$ws.Range("D43").Value = "'3.53"
$ws.Range("E43").Value = "  +6.05%  "
$ws.Range("D44").Value = "'5.10"
$ws.Range("E44").Value = "  +7.67%  "
$ws.Range("D45").Value = "'263.67"
$ws.Range("E45").Value = "  +11.28%  "
$ws.Range("D46").Value = "'0.0509"
$ws.Range("E46").Value = "  +4.74%  "
$ws.Range("E47").Value = "  +3.81%  "
$ws.Range("D48").Value = "'0.557"
$ws.Range("E48").Value = "  +3.42%  "
$ws.Range("D49").Value = "'0.0214"
$ws.Range("E49").Value = "  +5.61%  "
$ws.Range("D50").Value = "'0.378"
$ws.Range("E50").Value = "  +2.65%  "
$ws.Range("D51").Value = "'16.78"
$ws.Range("E51").Value = "  +5.43%  "
